# Automatic update of files.
#
# The source export rotated the three species-observation records that
# live in rows 4-6 of the "Artfynd" sheet: the record that was on row 4
# moved to row 6, the record on row 5 moved to row 4, and the record on
# row 6 moved to row 5 (i.e. a cyclic shift up by one row, with row 4's
# original data wrapping around to row 6). Re-create that by writing the
# new per-row values directly, cell by cell, instead of trying to "move"
# ranges (safer for the handful of cells that are blank in one row and
# populated in another).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 4  (becomes the old row 5 data)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 131136874
$ws.Range("B4").Value = 79244
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("J4").Value = "bålar"
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 788960
$ws.Range("R4").Value = 7131416
$ws.Range("AC4").ClearContents()
$ws.Range("AE4").Value = $false
# Writing a bare "" does not create a cell when none existed before (the
# engine drops empty-string writes to previously-absent cells), so force a
# real, present, empty *text* cell the same way Excel itself does: type an
# apostrophe (text-quote prefix), which stores as an empty string, then
# drop the quote-prefix style flag it leaves behind.
$ws.Range("AF4").Value = "'"
$ws.Range("AF4").Style = "Normal"

# ---------------------------------------------------------------------
# Row 5  (becomes the old row 6 data)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 131136941
$ws.Range("B5").Value = 83090
$ws.Range("E5").Value = 1312
$ws.Range("F5").Value = "Gammelgransskål"
$ws.Range("G5").Value = "Pseudographis pinicola"
$ws.Range("H5").Value = "(Nyl.) Rehm"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("Q5").Value = 788995
$ws.Range("R5").Value = 7131220
$ws.Range("AC5").Value = "på en gammal senvuxen gran"

# ---------------------------------------------------------------------
# Row 6  (becomes the old row 4 data)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 131136961
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").ClearContents()
$ws.Range("L6").Value = "'"
$ws.Range("L6").Style = "Normal"
$ws.Range("M6").Value = "'"
$ws.Range("M6").Style = "Normal"
$ws.Range("Q6").Value = 789068
$ws.Range("R6").Value = 7131245
$ws.Range("AC6").Value = "barksprätt på gammal gran"
$ws.Range("AE6").Value = $true
$ws.Range("AF6").ClearContents()
